# Natmi following Dr Hou advice
# Re-ran NATMI with an added "FAPs" sending cluster: the existing Ntng1-Lrrc4c
# row (sCs -> sCs) is recomputed, and a new row is added for FAPs -> sCs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs (sending) -> Ntng1/Lrrc4c -> sCs (target), recomputed metrics
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ntng1"
$ws.Range("C2").Value = "Lrrc4c"
$ws.Range("D2").Value = "sCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.050782
$ws.Range("H2").Value = 0.152346
$ws.Range("I2").Value = 0.02036584351830488
$ws.Range("J2").Value = 0.02036584351830488
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.652639666666667
$ws.Range("N2").Value = 7.957919
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1347063475526667
$ws.Range("R2").Value = 1.212357127974
$ws.Range("S2").Value = 0.02036584351830488
$ws.Range("T2").Value = 0.02036584351830488

# Row 3: sCs (sending) -> Ntng1/Lrrc4c -> sCs (target), new row
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Ntng1"
$ws.Range("C3").Value = "Lrrc4c"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.442706666666667
$ws.Range("H3").Value = 7.32812
$ws.Range("I3").Value = 0.9796341564816952
$ws.Range("J3").Value = 0.9796341564816952
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.652639666666667
$ws.Range("N3").Value = 7.957919
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 6.479620598031111
$ws.Range("R3").Value = 58.31658538228
$ws.Range("S3").Value = 0.9796341564816952
$ws.Range("T3").Value = 0.9796341564816952
